# Add "CO2 Cap Aux" and "CO2 Compress" unit-process rows to the
# "Unit Processes" sheet, inserted right after the existing
# "simple_CO2capture" row (row 41) and before "simple_CO2storage"
# (which, along with everything below it, shifts down by two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 42-43 (rows 42.. shift down to 44..),
# inheriting formatting from the row above as Excel normally does.
$ws.Rows("42:43").Insert()

# Row 42: duplicae_CO2capture / CO2 Cap Aux
$ws.Range("A42").Value2 = "duplicae_CO2capture"
$ws.Range("B42").Value2 = "CCS"
$ws.Range("C42").Value2 = "CO2 Capture & Compression"
$ws.Range("D42").Value2 = "CO2"
$ws.Range("E42").Value2 = "inflow"
$ws.Range("F42").Value2 = "data/steel/steel_simplified_var.xlsx"
$ws.Range("G42").Value2 = "CO2 Cap Aux"
$ws.Range("H42").Value2 = "data/steel/steel_simplified_calcs.xlsx"
$ws.Range("I42").Value2 = "CO2 Cap Aux"

# Row 43: simple_compression / CO2 Compress
$ws.Range("A43").Value2 = "simple_compression"
$ws.Range("B43").Value2 = "CCS"
$ws.Range("C43").Value2 = "CO2 Compression"
$ws.Range("D43").Value2 = "compressed CO2"
$ws.Range("E43").Value2 = "outflow"
$ws.Range("F43").Value2 = "data/steel/steel_simplified_var.xlsx"
$ws.Range("G43").Value2 = "CO2 Compress"
$ws.Range("H43").Value2 = "data/steel/steel_simplified_calcs.xlsx"
$ws.Range("I43").Value2 = "CO2 Compress"

# Match the saved selection state (bottom-right pane active cell moves
# from I47 to I43 now that the new rows sit at 42-43).
$ws.Range("I43").Select()
